# The "Diih" entry (row 47, a freshly-added / unverified member with all-zero
# war participation) is moved further down the roster, to just above
# "nivelador" (originally row 73), so it now lands on row 72. Every row
# between the old and new position (48..72) shifts up by one to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the row being relocated before it moves (Value2 gives a real
# snapshot of the underlying data, unlike Value which just echoes the
# property descriptor in this host).
$movedRow = $ws.Range("A47:H47").Value2

# Remove it from its original position; rows below shift up automatically.
$ws.Rows("47:47").Delete()

# Re-open a blank row just above "nivelador" (now at row 73 after the
# deletion-shift) and drop the captured data back in at its new home, row 72.
$ws.Rows("72:72").Insert()
$ws.Range("A72:H72").Value = $movedRow
